$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New KPI set values used by the added rows (49-52)
$kpiSetCap = "PoS 2019 - IC Cinema - CAP"
$kpiSetReg = "PoS 2019 - IC Cinema - REG"
$kpiName = "Combo other"
$atomicOldImage = "Coke and meal: Image"
$atomicNewImage = "Juice and meal: Image"
$atomicOldFood = "Coke and meal: Food"
$atomicNewFood = "Juice and meal: Food"

$formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D{0},"'', a.description=''",D{0},"'', a.display_text=''",D{0},"''  WHERE s.name=''",A{0},"'' AND k.display_text=''",B{0},"'' AND a.name=''",C{0},"'';")'

# Row 49: PoS 2019 - IC Cinema - CAP / Combo other / Image
$ws.Range("A49").Value = $kpiSetCap
$ws.Range("B49").Value = $kpiName
$ws.Range("C49").Value = $atomicOldImage
$ws.Range("D49").Value = $atomicNewImage
$ws.Range("E49").Formula = ($formula -f 49)

# Row 50: PoS 2019 - IC Cinema - CAP / Combo other / Food
$ws.Range("A50").Value = $kpiSetCap
$ws.Range("B50").Value = $kpiName
$ws.Range("C50").Value = $atomicOldFood
$ws.Range("D50").Value = $atomicNewFood
$ws.Range("E50").Formula = ($formula -f 50)

# Row 51: PoS 2019 - IC Cinema - REG / Combo other / Image
$ws.Range("A51").Value = $kpiSetReg
$ws.Range("B51").Value = $kpiName
$ws.Range("C51").Value = $atomicOldImage
$ws.Range("D51").Value = $atomicNewImage
$ws.Range("E51").Formula = ($formula -f 51)

# Row 52: PoS 2019 - IC Cinema - REG / Combo other / Food
$ws.Range("A52").Value = $kpiSetReg
$ws.Range("B52").Value = $kpiName
$ws.Range("C52").Value = $atomicOldFood
$ws.Range("D52").Value = $atomicNewFood
$ws.Range("E52").Formula = ($formula -f 52)

# Apply the same "dashed border" formatting used by the existing Atomic-Name-Old
# column (copied from an existing cell that already carries that style, e.g. C35)
$ws.Range("C35").Copy()
$ws.Range("C49").PasteSpecial(-4122)
$ws.Range("C35").Copy()
$ws.Range("C50").PasteSpecial(-4122)
$ws.Range("C35").Copy()
$ws.Range("C51").PasteSpecial(-4122)
$ws.Range("C35").Copy()
$ws.Range("C52").PasteSpecial(-4122)

# Apply the same "no border" formatting used elsewhere for the Atomic-Name-New
# column, then highlight it yellow to flag the changed value (new style)
$ws.Range("B38").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("B38").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("B38").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("B38").Copy()
$ws.Range("D52").PasteSpecial(-4122)
$ws.Range("D49").Interior.Color = 65535
$ws.Range("D50").Interior.Color = 65535
$ws.Range("D51").Interior.Color = 65535
$ws.Range("D52").Interior.Color = 65535

$ws.Application.CutCopyMode = 0

# Update selection to mirror where the new rows were added
$ws.Range("E49:E52").Select()
